$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$data = @(
    ,@(2, 'P4259', '{''eft:palgyi-lh-npo'', ''eft:dpal-gyi-lhun-po'', ''eft:ban-de-dpal-gyi-lhun-po''}')
    ,@(3, 'P8220', '{''eft:devacandra''}')
    ,@(4, 'P8212', '{''eft:devendraraksita''}')
    ,@(5, 'P8183', '{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}')
    ,@(6, 'P8219', '{''eft:visuddhasimha''}')
    ,@(7, 'P8151', '{''eft:gayadhara''}')
    ,@(8, 'P8245', '{''eft:buddhakaravarma''}')
    ,@(9, 'https://lod.dila.edu.tw/resource.php?id=A000089', '{''eft:siladharma''}')
    ,@(10, 'P00KG07267', '{''eft:sarvanyadeva'', ''eft:sarvajnadeva''}')
    ,@(11, 'P1242', '{''eft:g-wai-lodr-''}')
    ,@(12, 'P2551', '{''eft:blo-ldan-shes-rab''}')
    ,@(13, 'P8211', '{''eft:bidyakaraprabha'', ''eft:vidyakaraprabha''}')
    ,@(14, 'P4CZ15137', '{''eft:kumarakalasa''}')
    ,@(15, 'P8217', '{''eft:t-jnanagarbha'', ''eft:jnanagarbha''}')
    ,@(16, 'P8265', '{''eft:ratnaraksita''}')
    ,@(17, 'P8266', '{''eft:ch-nyi-tsultrim'', ''eft:lotsawa-band-dharmatasila'', ''eft:dharmatasila''}')
    ,@(18, 'P3458', '{''eft:g-lhets-''}')
    ,@(19, 'P3285', '{''eft:sakya-yesh-''}')
    ,@(20, 'P0TMPT007', '{''eft:rnam-par-mi-rtog-pa''}')
    ,@(21, 'P8221', '{''eft:g-ch-drup''}')
    ,@(22, 'P5788', '{''eft:th-nmi-sambhota''}')
    ,@(23, 'P3379', '{''eft:dipamkara-srijnana'', ''eft:dipamkarasrijnana''}')
    ,@(24, 'P2614', '{''eft:nyen-lotsawa-darma-drak''}')
    ,@(25, 'P3890', '{''eft:ch-kyi-sherab''}')
    ,@(26, 'P8261', '{''eft:munivarman'', ''eft:munivarma''}')
    ,@(27, 'P8267', '{''eft:vijayasila''}')
    ,@(28, 'P0TMP092', '{''eft:anandasri-s-''}')
    ,@(29, 'P4CZ15308', '{''eft:vairocanaraksita''}')
    ,@(30, 'P8278', '{''eft:gewai-lodr-'', ''eft:dge-ba-i-blo-gros''}')
    ,@(31, 'P8209', '{''eft:jinamitra-k-'', ''eft:jinamitra'', ''eft:dzi-na-mi-tra-k-''}')
    ,@(32, 'P0TMP098', '{''eft:jinavara''}')
    ,@(33, 'P3709', '{''eft:phakpa-sherab''}')
    ,@(34, 'P0TMP104', '{''eft:punyasambhava''}')
    ,@(35, 'P8263', '{''eft:leki-d-''}')
    ,@(36, 'P8276', '{''eft:wang-phab-zhwun-wang-phan-zhun-''}')
    ,@(37, 'P8268', '{''eft:buddhaprabha''}')
    ,@(38, 'P8277', '{''eft:rgya-mtsho-i-sde''}')
    ,@(39, 'P8205', '{''eft:zhang-yesh-d-'', ''eft:band-yesh-d-'', ''eft:yesh-de'', ''eft:yesh-d-'', ''eft:ye-shes-sde'', ''eft:band-yesh-de'', ''eft:yesh-d-ye-shes-sde-''}')
    ,@(40, 'P753', '{''eft:rin-chen-bzag-po'', ''eft:rinchen-zangpo'', ''eft:rin-chen-bzang-po''}')
    ,@(41, 'P3456', '{''eft:tshul-khrims-rgyal-ba'', ''eft:tsultrim-gyalwa''}')
    ,@(42, 'P8222', '{''eft:jnanasidhi'', ''eft:jnanasiddhi''}')
    ,@(43, 'P3214', '{''eft:danasila''}')
    ,@(44, '?', '{''eft:sherap-'', ''eft:sakyasena'', ''eft:vajrvisramitra''}')
    ,@(45, 'P4256', '{''eft:lotsawa-zangkyong-bzang-skyong-''}')
    ,@(46, 'P1KG8854', '{''eft:silendrabodhi'', ''eft:surendrabodhi'', ''eft:srilendrabodhi''}')
    ,@(47, 'P2956', '{''eft:krsnapandita''}')
    ,@(48, 'P4CZ16819', '{''eft:sakyaprabha''}')
    ,@(49, 'P5651', '{''eft:patsap-nyima-drak-'', ''eft:pa-tshab-nyi-ma-grags''}')
    ,@(50, 'P8269', '{''eft:dgon-gling-rma''}')
    ,@(51, 'P8271', '{''eft:kumararaksita''}')
    ,@(52, 'P4242', '{''eft:sherab-lekpa''}')
    ,@(53, 'P0TMP080', '{''eft:hwa-shang-zab-mo''}')
    ,@(54, 'P8280', '{''eft:subhasita''}')
    ,@(55, 'P8216', '{''eft:sakya-lodr-''}')
    ,@(56, 'P8249', '{''eft:pandita-dharmakara'', ''eft:dharmakara''}')
    ,@(57, 'P8273', '{''eft:rinchen-tso'', ''eft:rin-chen-tsho''}')
    ,@(58, 'P4255', '{''eft:t-jnanagarbha'', ''eft:ye-shes-snying-po'', ''eft:yesh-nyingpo''}')
    ,@(59, 'P6453', '{''eft:tsultrim-gyaltsen''}')
    ,@(60, 'P2557', '{''eft:-brom'', ''eft:-brom-rgyal-ba-i-byung-gnas''}')
    ,@(61, 'P8260', '{''eft:dpal-dbyangs''}')
    ,@(62, 'P8228', '{''eft:surendrabodhi''}')
    ,@(63, 'P2548', '{''eft:prajnavarma'', ''eft:prajnavarman''}')
    ,@(64, 'P8093', '{''eft:kamalagupta''}')
    ,@(65, 'P8182', '{''eft:ban-de-dpal-brtsegs'', ''eft:ska-ba-dpal-brtsegs'', ''eft:band-paltsek'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:dpal-brtsegs'', ''eft:paltsek''}')
    ,@(66, 'P8206', '{''eft:celu''}')
    ,@(67, 'P4263', '{''eft:gew-pal'', ''eft:dge-ba-dpal''}')
    ,@(68, 'P8213', '{''eft:t-vidyakarasimha'', ''eft:vidyakarasimha''}')
    ,@(69, 'P1321', '{''eft:shang-buchikpa''}')
    ,@(70, 'P2637', '{''eft:trakpa-gyaltsen''}')
    ,@(71, 'P0RK8', '{''eft:dharmapala''}')
    ,@(72, 'P4258', '{''eft:dpal-byor''}')
    ,@(73, 'P4CZ16780', '{''eft:manjusrigarbha''}')
    ,@(74, 'P8171', '{''eft:dharmasribhadra''}')
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
}
